$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Ref = "D2"; Val = "329.88" },
    @{ Ref = "E2"; Val = "6.96%" },
    @{ Ref = "G2"; Val = "7" },
    @{ Ref = "D3"; Val = "40.09" },
    @{ Ref = "E3"; Val = "7.19%" },
    @{ Ref = "G3"; Val = "7" },
    @{ Ref = "D4"; Val = "5.270" },
    @{ Ref = "E4"; Val = "1.91%" },
    @{ Ref = "G4"; Val = "7" },
    @{ Ref = "D5"; Val = "0.08082" },
    @{ Ref = "E5"; Val = "3.01%" },
    @{ Ref = "G5"; Val = "7" },
    @{ Ref = "D6"; Val = "4.517" },
    @{ Ref = "E6"; Val = "2.08%" },
    @{ Ref = "G6"; Val = "7" },
    @{ Ref = "D7"; Val = "8.649" },
    @{ Ref = "E7"; Val = "4.96%" },
    @{ Ref = "G7"; Val = "7" },
    @{ Ref = "D8"; Val = "1.937" },
    @{ Ref = "E8"; Val = "1.72%" },
    @{ Ref = "G8"; Val = "7" },
    @{ Ref = "E9"; Val = "-1.41%" },
    @{ Ref = "G9"; Val = "7" },
    @{ Ref = "E10"; Val = "0.03%" },
    @{ Ref = "G10"; Val = "7" },
    @{ Ref = "D11"; Val = "0.1369" },
    @{ Ref = "E11"; Val = "22.69%" },
    @{ Ref = "G11"; Val = "7" },
    @{ Ref = "D12"; Val = "0.1974" },
    @{ Ref = "E12"; Val = "1.11%" },
    @{ Ref = "G12"; Val = "7" },
    @{ Ref = "D13"; Val = "0.09092" },
    @{ Ref = "E13"; Val = "-0.09%" },
    @{ Ref = "G13"; Val = "7" },
    @{ Ref = "D14"; Val = "0.03505" },
    @{ Ref = "E14"; Val = "6.43%" },
    @{ Ref = "G14"; Val = "7" },
    @{ Ref = "D15"; Val = "0.09592" },
    @{ Ref = "E15"; Val = "-0.13%" },
    @{ Ref = "G15"; Val = "7" },
    @{ Ref = "D16"; Val = "0.001401" },
    @{ Ref = "E16"; Val = "1.41%" },
    @{ Ref = "G16"; Val = "7" },
    @{ Ref = "D17"; Val = "0.006355" },
    @{ Ref = "E17"; Val = "2.50%" },
    @{ Ref = "G17"; Val = "7" },
    @{ Ref = "D18"; Val = "3.365" },
    @{ Ref = "E18"; Val = "-6.88%" },
    @{ Ref = "G18"; Val = "7" },
    @{ Ref = "D19"; Val = "0.3518" },
    @{ Ref = "E19"; Val = "3.17%" },
    @{ Ref = "G19"; Val = "7" },
    @{ Ref = "D20"; Val = "6.492" },
    @{ Ref = "E20"; Val = "0.24%" },
    @{ Ref = "G20"; Val = "7" },
    @{ Ref = "E21"; Val = "2.59%" },
    @{ Ref = "G21"; Val = "7" },
    @{ Ref = "G22"; Val = "7" },
    @{ Ref = "D23"; Val = "0.04449" },
    @{ Ref = "E23"; Val = "0.89%" },
    @{ Ref = "G23"; Val = "7" },
    @{ Ref = "D24"; Val = "0.001223" },
    @{ Ref = "E24"; Val = "-0.75%" },
    @{ Ref = "G24"; Val = "7" },
    @{ Ref = "D25"; Val = "0.004324" },
    @{ Ref = "E25"; Val = "-5.35%" },
    @{ Ref = "G25"; Val = "7" },
    @{ Ref = "E26"; Val = "-0.80%" },
    @{ Ref = "G26"; Val = "7" },
    @{ Ref = "D27"; Val = "0.0003996" },
    @{ Ref = "E27"; Val = "0.10%" },
    @{ Ref = "G27"; Val = "7" },
    @{ Ref = "G28"; Val = "7" },
    @{ Ref = "G29"; Val = "7" },
    @{ Ref = "G30"; Val = "7" },
    @{ Ref = "G31"; Val = "7" },
    @{ Ref = "G32"; Val = "7" },
    @{ Ref = "G33"; Val = "7" },
    @{ Ref = "G34"; Val = "7" },
    @{ Ref = "G35"; Val = "7" },
    @{ Ref = "G36"; Val = "7" },
    @{ Ref = "G37"; Val = "7" },
    @{ Ref = "G38"; Val = "7" },
    @{ Ref = "D39"; Val = "0.02490" },
    @{ Ref = "E39"; Val = "12.99%" },
    @{ Ref = "G39"; Val = "7" },
    @{ Ref = "D40"; Val = "0.05200" },
    @{ Ref = "E40"; Val = "2.45%" },
    @{ Ref = "G40"; Val = "7" },
    @{ Ref = "D41"; Val = "0.007739" },
    @{ Ref = "E41"; Val = "3.78%" },
    @{ Ref = "G41"; Val = "7" },
    @{ Ref = "D42"; Val = "0.1426" },
    @{ Ref = "E42"; Val = "5.55%" },
    @{ Ref = "G42"; Val = "7" },
    @{ Ref = "D43"; Val = "0.009108" },
    @{ Ref = "E43"; Val = "4.03%" },
    @{ Ref = "G43"; Val = "7" },
    @{ Ref = "E44"; Val = "1.38%" },
    @{ Ref = "G44"; Val = "7" },
    @{ Ref = "D45"; Val = "0.009008" },
    @{ Ref = "E45"; Val = "4.35%" },
    @{ Ref = "G45"; Val = "7" },
    @{ Ref = "D46"; Val = "0.00006619" },
    @{ Ref = "E46"; Val = "1.03%" },
    @{ Ref = "G46"; Val = "7" },
    @{ Ref = "E47"; Val = "0.00%" },
    @{ Ref = "G47"; Val = "7" },
    @{ Ref = "D48"; Val = "0.003346" },
    @{ Ref = "E48"; Val = "16.92%" },
    @{ Ref = "G48"; Val = "7" },
    @{ Ref = "E49"; Val = "147.89%" },
    @{ Ref = "G49"; Val = "7" },
    @{ Ref = "E50"; Val = "0.00%" },
    @{ Ref = "G50"; Val = "7" },
    @{ Ref = "E51"; Val = "0.00%" },
    @{ Ref = "G51"; Val = "7" }
)

foreach ($u in $updates) {
    $ws.Range($u.Ref).Value = "'" + $u.Val
    $ws.Range($u.Ref).Style = "Normal"
}
